$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '27.496.60'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').Value = '1.825.39'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('E4').Value = '  -0.41%  '
Set-TextValue 'D5' '312.41'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  -0.65%  '
Set-TextValue 'D8' '0.3611'
$ws.Range('E8').Value = '  +0.80%  '
Set-TextValue 'D9' '0.07210'
$ws.Range('E9').Value = '  -1.13%  '
Set-TextValue 'D10' '0.8615'
$ws.Range('E10').Value = '  -1.12%  '
Set-TextValue 'D11' '20.64'
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('D12').Value = '1.805.41'
$ws.Range('E12').Value = '  -2.71%  '
Set-TextValue 'D13' '5.396'
$ws.Range('E13').Value = '  +1.18%  '
Set-TextValue 'D14' '6.484'
$ws.Range('E14').Value = '  -1.05%  '
Set-TextValue 'D15' '0.06927'
$ws.Range('E15').Value = '  -1.09%  '
$ws.Range('E16').Value = '  -0.55%  '
Set-TextValue 'D17' '80.55'
$ws.Range('E17').Value = '  +1.17%  '
Set-TextValue 'D18' '0.000008922'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('E19').Value = '  -0.36%  '
Set-TextValue 'D20' '15.38'
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('D21').Value = '27.486.27'
$ws.Range('E21').Value = '  -0.85%  '
Set-TextValue 'D22' '5.118'
$ws.Range('E22').Value = '  +2.39%  '
Set-TextValue 'D23' '10.90'
$ws.Range('E23').Value = '  +4.99%  '
$ws.Range('D24').Value = '2.049.11'
$ws.Range('E24').Value = '  -1.75%  '
$ws.Range('E25').Value = '  +0.22%  '
Set-TextValue 'D26' '154.93'
$ws.Range('E26').Value = '  -0.59%  '
$ws.Range('E27').Value = '  +0.90%  '
Set-TextValue 'D28' '5.166'
$ws.Range('E28').Value = '  -2.05%  '
Set-TextValue 'D29' '114.12'
$ws.Range('E29').Value = '  -5.24%  '
Set-TextValue 'D30' '1.793'
$ws.Range('E30').Value = '  -4.05%  '
Set-TextValue 'D31' '0.08872'
$ws.Range('E31').Value = '  -0.50%  '
Set-TextValue 'D32' '0.7500'
$ws.Range('E32').Value = '  -1.17%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D33' '4.545'
$ws.Range('E33').Value = '  +0.90%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D34' '2.975'
$ws.Range('E34').Value = '  +0.19%  '
Set-TextValue 'D35' '1.123'
$ws.Range('E35').Value = '  +0.15%  '
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('E37').Value = '  -1.06%  '
Set-TextValue 'D38' '0.05289'
$ws.Range('E38').Value = '  -2.55%  '
Set-TextValue 'D39' '0.01923'
$ws.Range('E39').Value = '  -0.37%  '
Set-TextValue 'D40' '2.792'
$ws.Range('E40').Value = '  -1.41%  '
Set-TextValue 'D41' '0.5077'
$ws.Range('E41').Value = '  +0.14%  '
Set-TextValue 'D42' '0.1656'
$ws.Range('E42').Value = '  -0.60%  '
Set-TextValue 'D43' '6.453'
$ws.Range('E43').Value = '  -2.24%  '
Set-TextValue 'D44' '8.359'
$ws.Range('E44').Value = '  -0.47%  '
Set-TextValue 'D45' '10.44'
$ws.Range('E45').Value = '  +0.53%  '
Set-TextValue 'D46' '106.23'
$ws.Range('E46').Value = '  +0.03%  '
Set-TextValue 'D48' '0.4686'
$ws.Range('E48').Value = '  +0.58%  '
Set-TextValue 'D49' '0.9999'
$ws.Range('E49').Value = '  -0.36%  '
Set-TextValue 'D50' '1.615'
$ws.Range('E50').Value = '  -0.84%  '
Set-TextValue 'D51' '63.82'
$ws.Range('E51').Value = '  -1.10%  '
